$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.696.06"
$ws.Range("E2").Value = "  +5.13%  "

$ws.Range("D3").Value = "2.584.66"
$ws.Range("E3").Value = "  +5.60%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.37%  "

$ws.Range("D5").Value = "'590.31"
$ws.Range("E5").Value = "  +3.58%  "

$ws.Range("D6").Value = "'155.14"
$ws.Range("E6").Value = "  +6.55%  "

$ws.Range("E7").Value = "  -0.35%  "

$ws.Range("D8").Value = "'0.547"
$ws.Range("E8").Value = "  +3.29%  "

$ws.Range("D9").Value = "2.609.89"
$ws.Range("E9").Value = "  +6.67%  "

$ws.Range("D10").Value = "'0.115"
$ws.Range("E10").Value = "  +4.10%  "

$ws.Range("E11").Value = "  -1.52%  "

$ws.Range("E12").Value = "  +4.71%  "

$ws.Range("D13").Value = "'5.31"
$ws.Range("E13").Value = "  +2.18%  "

$ws.Range("D14").Value = "'29.29"
$ws.Range("E14").Value = "  +2.60%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.063.11"
$ws.Range("E15").Value = "  +6.15%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000181"
$ws.Range("E16").Value = "  +4.51%  "

$ws.Range("D17").Value = "65.723.12"
$ws.Range("E17").Value = "  +5.37%  "

$ws.Range("D18").Value = "2.614.65"
$ws.Range("E18").Value = "  +6.67%  "

$ws.Range("D19").Value = "'8.21"
$ws.Range("E19").Value = "  +8.32%  "

$ws.Range("D20").Value = "'11.23"
$ws.Range("E20").Value = "  +4.82%  "

$ws.Range("D21").Value = "'356.34"
$ws.Range("E21").Value = "  +11.26%  "

$ws.Range("E22").Value = "  +4.65%  "

$ws.Range("D23").Value = "'2.24"
$ws.Range("E23").Value = "  +2.39%  "

$ws.Range("E24").Value = "  -0.29%  "

$ws.Range("D25").Value = "'9.98"
$ws.Range("E25").Value = "  +0.83%  "

$ws.Range("D26").Value = "'66.26"
$ws.Range("E26").Value = "  +2.23%  "

$ws.Range("D27").Value = "'632.94"
$ws.Range("E27").Value = "  -1.45%  "

$ws.Range("D28").Value = "'0.0000105"
$ws.Range("E28").Value = "  +10.60%  "

$ws.Range("D29").Value = "2.714.78"
$ws.Range("E29").Value = "  +6.04%  "

$ws.Range("D30").Value = "'1.50"
$ws.Range("E30").Value = "  +7.01%  "

$ws.Range("D31").Value = "'0.996"
$ws.Range("E31").Value = "  +1.37%  "

$ws.Range("D32").Value = "'8.23"
$ws.Range("E32").Value = "  +5.36%  "

$ws.Range("D33").Value = "'1.90"
$ws.Range("E33").Value = "  +5.70%  "

$ws.Range("D34").Value = "'0.137"
$ws.Range("E34").Value = "  +3.99%  "

$ws.Range("D35").Value = "'1.63"
$ws.Range("E35").Value = "  +9.64%  "

$ws.Range("D36").Value = "'0.994"
$ws.Range("E36").Value = "  -0.28%  "

$ws.Range("D37").Value = "'4.97"
$ws.Range("E37").Value = "  +7.33%  "

$ws.Range("D38").Value = "'5.64"
$ws.Range("E38").Value = "  +7.07%  "

$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").Value = "'19.38"
$ws.Range("E39").Value = "  +5.14%  "

$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "'2.88"
$ws.Range("E40").Value = "  +6.95%  "

$ws.Range("D41").Value = "'155.73"
$ws.Range("E41").Value = "  +3.58%  "

$ws.Range("E42").Value = "  +2.96%  "

$ws.Range("D43").Value = "'1.82"
$ws.Range("E43").Value = "  +6.68%  "

$ws.Range("D44").Value = "0.0₆0326"
$ws.Range("E44").Value = "  +5.83%  "

$ws.Range("D45").Value = "'42.05"
$ws.Range("E45").Value = "  +0.93%  "

$ws.Range("D46").Value = "'163.52"
$ws.Range("E46").Value = "  +7.66%  "

$ws.Range("E47").Value = "  -0.12%  "

$ws.Range("D48").Value = "'16.20"
$ws.Range("E48").Value = "  +5.62%  "

$ws.Range("D49").Value = "'3.76"
$ws.Range("E49").Value = "  +6.74%  "

$ws.Range("D50").Value = "'21.75"
$ws.Range("E50").Value = "  +8.63%  "

$ws.Range("D51").Value = "'0.635"
$ws.Range("E51").Value = "  +5.59%  "
